# Create separate columns for system 4 and system 6 heating coil types.
#
# Originally column J ("heating_coil_type_sys4and6") held a single fuel-type
# default shared by both system 4 and system 6. This change splits it into
# two columns: J becomes "heating_coil_type_sys4" (keeps the original values)
# and a brand-new column K becomes "heating_coil_type_sys6" (defaults to the
# hot-water-coil equivalent for gas-heated rows, and the same electric default
# for electrically-heated rows). The former column K ("fan_type") and
# everything to its right shifts one column to the right (now column L, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank column at K; this shifts the former K (fan_type) and
# everything after it one column to the right, and - just like Excel's
# native "Insert Copied/Cut Cells" behaviour - the new column inherits the
# formatting of the column immediately to its left (J) on a per-row basis.
$ws.Columns.Item(11).Insert()

# Rename the headers.
$ws.Cells.Item(1, 10).Value = "heating_coil_type_sys4"
$ws.Cells.Item(1, 11).Value = "heating_coil_type_sys6"

$lastRow = $ws.UsedRange.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $sys4 = $ws.Cells.Item($r, 10).Value()

    if ($sys4 -eq "Gas") {
        $sys6 = "Hot Water"
    } else {
        $sys6 = $sys4
    }

    $ws.Cells.Item($r, 11).Value = $sys6
}
